$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure updated price/volume cells keep their original plain-text representation
# (avoids Excel auto-converting numeric-looking strings into floating point numbers
# or percentages, which would corrupt precision / trailing zeros).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "306.47"
Set-TextValue "E2" "-0.60%"
Set-TextValue "D3" "39.12"
Set-TextValue "E3" "7.74%"
Set-TextValue "D4" "5.093"
Set-TextValue "E4" "0.79%"
Set-TextValue "D5" "0.08053"
Set-TextValue "E5" "-0.67%"
Set-TextValue "D6" "1.957"
Set-TextValue "E6" "-5.40%"
Set-TextValue "D7" "4.194"
Set-TextValue "E7" "0.92%"
Set-TextValue "D8" "8.005"
Set-TextValue "E8" "1.93%"
Set-TextValue "D9" "0.9319"
Set-TextValue "E9" "0.18%"
Set-TextValue "D10" "0.1443"
Set-TextValue "E10" "1.85%"
Set-TextValue "D11" "0.1933"
Set-TextValue "E11" "-0.01%"
Set-TextValue "D12" "0.09126"
Set-TextValue "E12" "0.61%"
Set-TextValue "D13" "0.03506"
Set-TextValue "E13" "1.56%"
Set-TextValue "D14" "0.09784"
Set-TextValue "E14" "-1.36%"
Set-TextValue "D15" "0.001393"
Set-TextValue "E15" "-0.91%"
Set-TextValue "D16" "0.006053"
Set-TextValue "E16" "-4.14%"
Set-TextValue "E17" "-1.33%"
Set-TextValue "E18" "-1.49%"
Set-TextValue "E19" "-0.63%"
Set-TextValue "D20" "0.1347"
Set-TextValue "E20" "5.04%"
Set-TextValue "D21" "4.569"
Set-TextValue "E21" "-4.99%"
Set-TextValue "D23" "0.04371"
Set-TextValue "E23" "-0.13%"
Set-TextValue "D25" "0.004275"
Set-TextValue "E25" "-13.10%"
Set-TextValue "E26" "0.14%"
Set-TextValue "D39" "0.02035"
Set-TextValue "E39" "0.28%"
Set-TextValue "D40" "0.05057"
Set-TextValue "E40" "-2.07%"
Set-TextValue "D41" "0.007444"
Set-TextValue "E41" "-0.73%"
Set-TextValue "D42" "0.01026"
Set-TextValue "E42" "2.17%"
Set-TextValue "D43" "0.1347"
Set-TextValue "E43" "-1.78%"
Set-TextValue "E44" "-2.17%"
Set-TextValue "D45" "0.009118"
Set-TextValue "E45" "-8.61%"
Set-TextValue "D46" "0.00006202"
Set-TextValue "E46" "-1.29%"
Set-TextValue "E47" "0.20%"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.20%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.20%"
